$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D3 (existing "28.01: Oppgaveseminar..." text) with appended sentence
$ws.Range("D3").Value = "28.01: Oppgaveseminar på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09). Se \@ref(seminar) for oppgaver."

# Add new task set for week 5 (row 4) and week 6 (row 5)
$ws.Range("C4").Value = "02.02: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("D4").Value = "04.02: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("C5").Value = "09.02: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("D5").Value = "11.02: Oppgaveseminar på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09). Se \@ref(seminar) for oppgaver."

# Update selection to match the committed state
$ws.Range("F17").Select()
